$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 405
$ws.Range("F5").Value = 1102
$ws.Range("F8").Value = 820
$ws.Range("F9").Value = 1584
$ws.Range("F10").Value = 6004
$ws.Range("G10").Value = 19.9
$ws.Range("F11").Value = 107
$ws.Range("F12").Value = 1714
$ws.Range("F13").Value = 431
$ws.Range("F14").Value = 5793
$ws.Range("F15").Value = 5793
$ws.Range("F20").Value = 1620
$ws.Range("F24").Value = 1305
$ws.Range("F25").Value = 710
$ws.Range("F26").Value = 217
$ws.Range("F29").Value = 23
$ws.Range("F30").Value = 74

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 297
$ws.Range("F5").Value = 149
$ws.Range("F8").Value = 364

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2208
$ws.Range("F5").Value = 145

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2208
$ws.Range("F6").Value = 405
$ws.Range("F7").Value = 1102
$ws.Range("F11").Value = 297
$ws.Range("F12").Value = 820
$ws.Range("F13").Value = 145
$ws.Range("F14").Value = 1584
$ws.Range("F15").Value = 6004
$ws.Range("G15").Value = 19.9
$ws.Range("F16").Value = 107
$ws.Range("F17").Value = 364
$ws.Range("F18").Value = 1714
$ws.Range("F21").Value = 431
$ws.Range("F24").Value = 5793
$ws.Range("F25").Value = 5793
$ws.Range("F30").Value = 1620
$ws.Range("F34").Value = 1305
$ws.Range("F35").Value = 710
$ws.Range("F36").Value = 217
$ws.Range("F44").Value = 74
